$d = $word.ActiveDocument

# --- Edit 1 ---
# Insert a brand-new bullet paragraph before the "Navigate to ... AutomationTools ... automation.py"
# paragraph, describing creating a new git branch first.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*automation.py*") {
        $p.Range.InsertBefore("Create a new branch for the episode off of the master branch, you can use gitkraken to make this easily`r")
        break
    }
}

# --- Edit 2 ---
# Change "Push to git" -> "Push to the new git branch", then add a new bullet paragraph
# after it about creating a pull request into master.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Push to git*") {
        $p.Range.Text = "Push to the new git branch"
        $p.Range.InsertAfter("`rCreate a pull request into master")
        break
    }
}
